# Generate Report for Handback
# Adds a new handback entry (0617a83e-cf74-4a24-8951-5e21e122a888) as row 4
# on the "Overview", "zh-cn" and "de-de" sheets, expanding each sheet's
# table (ListObject) to include the new row.

$wb = $excel.ActiveWorkbook

$blueColor = 15570276   # matches the workbook's existing HyperLink font color FF6495ED
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)
[void]$tblOverview.ListRows.Add()

$wsOverview.Range("A4").Value = "0617a83e-cf74-4a24-8951-5e21e122a888.md"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/923533a5f697198aab851c6136e03aea0e968096/e2e/0617a83e-cf74-4a24-8951-5e21e122a888.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "e2e\0617a83e-cf74-4a24-8951-5e21e122a888.md"
) | Out-Null
$wsOverview.Range("B4").Font.Underline = $true
$wsOverview.Range("B4").Font.Color = $blueColor

$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"

$wsOverview.Range("G4").Value = "2016-11-08 23:42:29"
$wsOverview.Range("G4").NumberFormat = $dateFormat

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$tblZhCn = $wsZhCn.ListObjects.Item(1)
[void]$tblZhCn.ListRows.Add()

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/923533a5f697198aab851c6136e03aea0e968096/e2e/0617a83e-cf74-4a24-8951-5e21e122a888.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "0617a83e-cf74-4a24-8951-5e21e122a888.md"
) | Out-Null
$wsZhCn.Range("A4").Font.Underline = $true
$wsZhCn.Range("A4").Font.Color = $blueColor

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "0617a83e-cf74-4a24-8951-5e21e122a888.81f540de6c8b41bf746b7c23e0073621270a5d04.zh-cn.xlf"

$wsZhCn.Range("H4").Value = "2016-11-08 23:42:15"
$wsZhCn.Range("H4").NumberFormat = $dateFormat

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/81f540de6c8b41bf746b7c23e0073621270a5d04/e2e/0617a83e-cf74-4a24-8951-5e21e122a888.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "0617a83e-cf74-4a24-8951-5e21e122a888.md"
) | Out-Null
$wsZhCn.Range("I4").Font.Underline = $true
$wsZhCn.Range("I4").Font.Color = $blueColor

$wsZhCn.Range("J4").Value = "0617a83e-cf74-4a24-8951-5e21e122a888.81f540de6c8b41bf746b7c23e0073621270a5d04.zh-cn.xlf"

$wsZhCn.Range("K4").Value = "2016-11-08 23:43:07"
$wsZhCn.Range("K4").NumberFormat = $dateFormat

$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$tblDeDe = $wsDeDe.ListObjects.Item(1)
[void]$tblDeDe.ListRows.Add()

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/923533a5f697198aab851c6136e03aea0e968096/e2e/0617a83e-cf74-4a24-8951-5e21e122a888.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "0617a83e-cf74-4a24-8951-5e21e122a888.md"
) | Out-Null
$wsDeDe.Range("A4").Font.Underline = $true
$wsDeDe.Range("A4").Font.Color = $blueColor

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "0617a83e-cf74-4a24-8951-5e21e122a888.81f540de6c8b41bf746b7c23e0073621270a5d04.de-de.xlf"

$wsDeDe.Range("H4").Value = "2016-11-08 23:42:29"
$wsDeDe.Range("H4").NumberFormat = $dateFormat

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/81f540de6c8b41bf746b7c23e0073621270a5d04/e2e/0617a83e-cf74-4a24-8951-5e21e122a888.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "0617a83e-cf74-4a24-8951-5e21e122a888.md"
) | Out-Null
$wsDeDe.Range("I4").Font.Underline = $true
$wsDeDe.Range("I4").Font.Color = $blueColor

$wsDeDe.Range("J4").Value = "0617a83e-cf74-4a24-8951-5e21e122a888.81f540de6c8b41bf746b7c23e0073621270a5d04.de-de.xlf"

$wsDeDe.Range("K4").Value = "2016-11-08 23:43:26"
$wsDeDe.Range("K4").NumberFormat = $dateFormat

$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"
